$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 1. "...involved the use of the following methods: Fictional Narrative,
#     Personas and Extreme Personas."
#  -> "...involved the use of the following methods: Personas and Extreme
#     Personas."
# ------------------------------------------------------------------
$d.Content.Find.Execute("Fictional Narrative, Personas and Extreme Personas.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Personas and Extreme Personas.", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Remove the "Fictional Narrative" heading paragraph together with the
#    blank paragraph that follows it (the section itself was dropped).
# ------------------------------------------------------------------
$pFictional = $d.Paragraphs(24)
$pBlankAfterFictional = $d.Paragraphs(25)
$d.Range($pFictional.Range.Start, $pBlankAfterFictional.Range.End).Delete() | Out-Null

# ------------------------------------------------------------------
# 3. The "_GoBack" bookmark used to sit near the end of the "gamification"
#    paragraph; it now sits in the (now empty) paragraph left behind by the
#    deletion above. Re-anchoring it here (by name) moves it from its old
#    location.
# ------------------------------------------------------------------
$pNowEmpty = $d.Paragraphs(23)
$d.Bookmarks.Add("_GoBack", $pNowEmpty.Range) | Out-Null

# ------------------------------------------------------------------
# 4. Insert the new "Personas is a very effective..." paragraph (plus its
#    trailing blank line) right before the "Extreme Personas" heading.
# ------------------------------------------------------------------
$pExtremePersonas = $d.Paragraphs(26)
$insertBeforeExtreme = $d.Range($pExtremePersonas.Range.Start, $pExtremePersonas.Range.Start)
$personasXml = "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">Personas is a very effective method that is widely used </w:t></w:r>" +
    "<w:r><w:t>in the design process to instigate some initial ideas.</w:t></w:r>" +
    "</w:p><w:p $wns/>"
$insertBeforeExtreme.InsertXML($personasXml)

# ------------------------------------------------------------------
# 5. Trailing "    " (four spaces) after "...minimize injuries." becomes a
#    single space. This also clears out the old "_GoBack" bookmark that
#    used to live in that run (already re-anchored in step 3).
# ------------------------------------------------------------------
$d.Content.Find.Execute("to minimize injuries.    ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "to minimize injuries. ", 2) | Out-Null

# ------------------------------------------------------------------
# 6. Insert the new blank line + "After some brainstorming..." paragraph
#    right before the pre-existing blank paragraph that precedes
#    "Final Design Solution". The fragment's final (non-empty) paragraph
#    is absorbed into the paragraph immediately following the insertion
#    point, so anchoring on that pre-existing blank paragraph lets it pick
#    up our new text while leaving "Final Design Solution" itself
#    untouched.
# ------------------------------------------------------------------
$pBlankBeforeFinalDesign = $d.Paragraphs(35)
$insertBeforeBlank = $d.Range($pBlankBeforeFinalDesign.Range.Start, $pBlankBeforeFinalDesign.Range.Start)
$afterXml = "<w:p $wns/><w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">After some brainstorming for gamification ideas to apply to the product, I concluded that adding some type of scoring system to the technique checker. These scores could then be uploaded and a ladder could be compiled, showing the top scoring </w:t></w:r>" +
    "<w:r><w:t>users. This</w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`"> would </w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`">certainly </w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`">add further motivation to </w:t></w:r>" +
    "<w:r><w:t>system and also attract geeks and gamers to possibly start using the system and commence working out themselves as now there is a gaming and competition aspect to the system.</w:t></w:r>" +
    "</w:p>"
$insertBeforeBlank.InsertXML($afterXml)
